$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 2038323.67
$ws.Range("C7").Value = -54.87871599249057
$ws.Range("D7").Value = 1951
$ws.Range("E7").Value = 1951
$ws.Range("F7").Value = 1044.758416196822
$ws.Range("G7").Value = 7.819285516662711
